$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("A1").Value = "Leader_Expr"
$ws.Range("B1").Value = "Follower_Expr"
$ws.Range("A2").Value = "x"
$ws.Range("B2").Value = "y"

$ws = $wb.Worksheets.Item(2)
$ws.Range("A1").Value = "Expression"
$ws.Range("B1").Value = "Function_Evaluation"
$ws.Range("C1").Value = "Restriction_Set_Type"
$ws.Range("D1").Value = "MIU_value"
$ws.Range("A2").Value = "2.09 - x"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "-3.09"
$ws.Range("C2").Value = "J_0_g"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0.86"
$ws.Range("A3").Value = "-2.09 + x"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "1.0899999999999999"
$ws.Range("C3").Value = "J_0_g"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "0.62"
$ws.Range("A4").Value = "41.02289999999999 + x - y - 9(x^2)"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "-40.02289999999999"
$ws.Range("C4").Value = "J_0_g"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.58"

$ws = $wb.Worksheets.Item(3)
$ws.Range("A1").Value = "Expression"
$ws.Range("B1").Value = "Function_Evaluation"
$ws.Range("C1").Value = "Restriction_Set_Type"
$ws.Range("D1").Value = "Lambda_value"
$ws.Range("E1").Value = "Beta_value"
$ws.Range("F1").Value = "Gamma_value"
$ws.Range("A2").Value = "12.69799662447257 - 9.383578059071729y + (-0.5 + x)*(y^2)"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "-12.69799662447257"
$ws.Range("C2").Value = "J_0_L0_v"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0.69"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "6.4"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "6.4"
$ws.Range("A3").Value = "5.206 - 1.37y"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "-6.206"
$ws.Range("C3").Value = "J_0_LP_v"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "0.65"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.8"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "1.2"
$ws.Range("A4").Value = "-2.4950632911392407 + 0.6329113924050633y"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "1.4050632911392404"
$ws.Range("C4").Value = "J_Ne_L0_v"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.32"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4.5"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "1.5"

$ws = $wb.Worksheets.Item(4)
$ws.Range("A1").Value = "x"
$ws.Range("B1").Value = "y"
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2.09"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "3.8"

$ws = $wb.Worksheets.Item(5)
$ws.Range("A1").Value = "vec_bf"
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "-2.1753227848101258"

$ws = $wb.Worksheets.Item(6)
$ws.Range("A1").Value = "vec_BF"
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "-84.36468"
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "-18.175055687763706"

$ws = $wb.Worksheets.Item(7)
$ws.Range("A1").Value = "vec_alpha"

$ws = $wb.Worksheets.Item(7)
$ws.Range("A2").Value = 2.37
